$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 139 (existing rows 139-157 shift down to 142-160).
$ws.Range("A139:T141").EntireRow.Insert()

# --- New row 139 ---
$ws.Range("A139").Value = 1
$ws.Range("B139").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C139").Value = "Arica y Parinacota"
$ws.Range("D139").Value = 45127
$ws.Range("E139").Value = 15
$ws.Range("F139").Value = "Fruta"
$ws.Range("G139").Value = 100102
$ws.Range("H139").Value = "Cítricos"
$ws.Range("I139").Value = 100102005
$ws.Range("J139").Value = "Naranja"
$ws.Range("K139").Value = "Fukumoto"
$ws.Range("L139").Value = "Primera"
$ws.Range("M139").Value = 250
$ws.Range("N139").Value = 900
$ws.Range("O139").Value = 950
$ws.Range("P139").Value = 925
$ws.Range("Q139").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R139").Value = "Región de O'Higgins"
$ws.Range("S139").Value = 925
$ws.Range("T139").Value = 1

# --- New row 140 ---
$ws.Range("A140").Value = 1
$ws.Range("B140").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C140").Value = "Arica y Parinacota"
$ws.Range("D140").Value = 45127
$ws.Range("E140").Value = 15
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100102
$ws.Range("H140").Value = "Cítricos"
$ws.Range("I140").Value = 100102005
$ws.Range("J140").Value = "Naranja"
$ws.Range("K140").Value = "Fukumoto"
$ws.Range("L140").Value = "Segunda"
$ws.Range("M140").Value = 300
$ws.Range("N140").Value = 750
$ws.Range("O140").Value = 850
$ws.Range("P140").Value = 800
$ws.Range("Q140").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R140").Value = "Región de O'Higgins"
$ws.Range("S140").Value = 800
$ws.Range("T140").Value = 1

# --- New row 141 ---
$ws.Range("A141").Value = 1
$ws.Range("B141").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C141").Value = "Arica y Parinacota"
$ws.Range("D141").Value = 45127
$ws.Range("E141").Value = 15
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100102
$ws.Range("H141").Value = "Cítricos"
$ws.Range("I141").Value = 100102005
$ws.Range("J141").Value = "Naranja"
$ws.Range("K141").Value = "Fukumoto"
$ws.Range("L141").Value = "Tercera"
$ws.Range("M141").Value = 300
$ws.Range("N141").Value = 700
$ws.Range("O141").Value = 750
$ws.Range("P141").Value = 725
$ws.Range("Q141").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R141").Value = "Región de O'Higgins"
$ws.Range("S141").Value = 725
$ws.Range("T141").Value = 1
